# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E69) listed the 54 monthly periods from
# 1802 (Feb-2018) through 2207 (Jul-2022) in ascending order. The sheet was
# updated to list them in descending order (2207 down to 1802) instead, and
# the "Valor Mora" figures in column F follow the rows, so the 66667 value
# that used to sit next to period 2207 (row 69) now sits next to period 2207
# at the top (row 16), and the 80000 value that used to sit next to 1802
# (row 16) now sits next to 1802 at the bottom (row 69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the list of periods 1802..2207 (YYMM, Feb 2018 through Jul 2022)
# in the same ascending order they originally appear in rows 16..69.
$periods = @()
for ($y = 18; $y -le 22; $y++) {
    for ($m = 1; $m -le 12; $m++) {
        if ($y -eq 18 -and $m -lt 2) { continue }
        if ($y -eq 22 -and $m -gt 7) { continue }
        $periods += ("{0:D2}{1:D2}" -f $y, $m)
    }
}

# Reverse it so the newest period (2207) comes first.
$periodsDesc = $periods[($periods.Count - 1)..0]

$firstRow = 16
for ($i = 0; $i -lt $periodsDesc.Count; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periodsDesc[$i]
}

# Swap the two "Valor Mora" amounts that no longer line up with their
# original periods now that the period column has been reversed.
$lastRow = $firstRow + $periodsDesc.Count - 1
$valTop = $ws.Cells.Item($firstRow, 6).Value2
$valBottom = $ws.Cells.Item($lastRow, 6).Value2
$ws.Cells.Item($firstRow, 6).Value = $valBottom
$ws.Cells.Item($lastRow, 6).Value = $valTop

$wb.Save()
